# Fruta / hortaliza, semanal
# Insert a new record row at row 87, pushing existing rows 87-90 down to 88-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 87 (shifts old rows 87-90 down to 88-91)
$ws.Rows.Item(87).Insert()

# Populate the new row 87 with the new record's data
$ws.Range("A87").Value = 7
$ws.Range("B87").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C87").Value = "Ñuble"
$ws.Range("D87").Value = 44783
$ws.Range("D87").NumberFormat = $ws.Range("D88").NumberFormat
$ws.Range("E87").Value = 16
$ws.Range("F87").Value = 100112031
$ws.Range("G87").Value = "Poroto verde"
$ws.Range("H87").Value = "Magnum"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 50
$ws.Range("K87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("M87").Value = 30000
$ws.Range("N87").Value = "$/malla 25 kilos"
$ws.Range("O87").Value = "Perú"
$ws.Range("P87").Value = 1200
$ws.Range("Q87").Value = 25
$ws.Range("R87").Value = "Hortaliza"
